# ==========================================================================
# Applies the "Changed automatic fill in for Walter-Lieth, added Samoa data,
# and added files to Getting Started and Learn More pages" edit.
# ==========================================================================

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Getting Started" summary sheet
$ws2 = $wb.Worksheets.Item(2)   # "Learn More" detail sheet

# --------------------------------------------------------------------------
# Sheet1 ("Getting Started"): remove the Hawaii rows, add a Samoa (climate)
# row right after the New Zealand (climate) row / before the Community rows.
# --------------------------------------------------------------------------

# Before:
#  1 DataType  | Area
#  2 Climate   | Hawaii (climate)
#  3 Climate   | Logan, UT (climate)
#  4 Climate   | New Zealand (climate)
#  5 Community | Hawaii (community)
#  6 Community | Logan, UT (community)
#  7 Community | New Zealand (community)
#
# After:
#  1 DataType  | Area
#  2 Climate   | Logan, UT (climate)
#  3 Climate   | New Zealand (climate)
#  4 Climate   | Samoa (climate)
#  5 Community | Logan, UT (community)
#  6 Community | New Zealand (community)

$ws1.Rows.Item(2).Delete()            # remove "Climate | Hawaii (climate)"
$ws1.Rows.Item(4).Delete()            # remove "Community | Hawaii (community)"
$ws1.Rows.Item(4).Insert()            # make room for the new Samoa row

$ws1.Cells.Item(4, 1).Value2 = "Climate"
$ws1.Cells.Item(4, 2).Value2 = "Samoa (climate)"

[void]$ws1.Range("B4").Select()

# --------------------------------------------------------------------------
# Sheet2 ("Learn More"): drop the Hawaii climate rows, turn the old Hawaii
# community rows into the new Samoa climate rows, and add a description row
# under both the Logan (FIA) and New Zealand (NVS) community entries.
# --------------------------------------------------------------------------

# Remove the two "Hawaii (climate)" rows (rows 2-3).
$ws2.Range("A2:A3").EntireRow.Delete()

# After the delete above the table looks like:
#  2-12  Logan / New Zealand climate rows (unchanged)
#  13-14 Hawaii (community)            -> will become Samoa (climate)
#  15    Logan, UT (community) eBird   (unchanged)
#  16    Logan, UT (community) plants  -> becomes 2 rows (FIA dataset + description)
#  17    New Zealand (community) eBird (unchanged)
#  18    New Zealand (community) plants-> becomes 2 rows (NVS dataset + description)

# Make room for the extra description rows.
$ws2.Rows.Item(17).Insert()           # room under the Logan FIA row
$ws2.Rows.Item(20).Insert()           # room under the New Zealand NVS row

# ---- Row 13: Samoa (climate) / Samoa ----
$ws2.Cells.Item(13, 1).Value2 = "Samoa (climate)"
$ws2.Cells.Item(13, 2).Value2 = "Samoa"
$ws2.Cells.Item(13, 3).Value2 = "Samoa_climate.xlsx"
$ws2.Cells.Item(13, 4).Value2 = -13.759
$ws2.Cells.Item(13, 5).Value2 = -172.1046
$ws2.Cells.Item(13, 6).Value2 = 23

# ---- Row 14: Samoa (climate) / XXXXXXX placeholder ----
$ws2.Cells.Item(14, 1).Value2 = "Samoa (climate)"
$ws2.Cells.Item(14, 2).Value2 = "XXXXXXX"
$ws2.Cells.Item(14, 3).Value2 = "XXXXXXX"
$ws2.Range("D14:F14").ClearContents()

# ---- Row 15: Logan, UT (community) / UT eBird Observation Dataset (unchanged) ----
$ws2.Cells.Item(15, 1).Value2 = "Logan, UT (community)"
$ws2.Cells.Item(15, 2).Value2 = "UT eBird Observation Dataset"
$ws2.Cells.Item(15, 3).Value2 = "Logan_eBird.xlsx"
$ws2.Cells.Item(15, 4).Value2 = 41.57
$ws2.Cells.Item(15, 5).Value2 = -111.7

# ---- Row 16: Logan, UT (community) / Forest Inventory and Analysis Dataset ----
$ws2.Cells.Item(16, 1).Value2 = "Logan, UT (community)"
$ws2.Cells.Item(16, 2).Value2 = "Forest Inventory and Analysis Dataset"
$ws2.Cells.Item(16, 3).Value2 = "Utah_Vegetation.xlsx"
$ws2.Range("D16:F16").ClearContents()

# ---- Row 17: Logan, UT (community) / FIA Data Description ----
$ws2.Cells.Item(17, 1).Value2 = "Logan, UT (community)"
$ws2.Cells.Item(17, 2).Value2 = "FIA Data Description"
$ws2.Cells.Item(17, 3).Value2 = "Utah_Vegetation_Description.xlsx"
$ws2.Range("D17:F17").ClearContents()

# ---- Row 18: New Zealand (community) / NZ eBird Observation Dataset (unchanged) ----
$ws2.Cells.Item(18, 1).Value2 = "New Zealand (community)"
$ws2.Cells.Item(18, 2).Value2 = "NZ eBird Observation Dataset"
$ws2.Cells.Item(18, 3).Value2 = "New_Zealand_eBird.xlsx"
$ws2.Cells.Item(18, 4).Value2 = -44.2
$ws2.Cells.Item(18, 5).Value2 = 170.5

# ---- Row 19: New Zealand (community) / National Vegetation Survey ----
$ws2.Cells.Item(19, 1).Value2 = "New Zealand (community)"
$ws2.Cells.Item(19, 2).Value2 = "National Vegetation Survey"
$ws2.Cells.Item(19, 3).Value2 = "NZ_Vegetation.xlsx"
$ws2.Range("D19:F19").ClearContents()

# ---- Row 20: New Zealand (community) / NVS Data Description ----
$ws2.Cells.Item(20, 1).Value2 = "New Zealand (community)"
$ws2.Cells.Item(20, 2).Value2 = "NVS Data Description"
$ws2.Cells.Item(20, 3).Value2 = "NZ_Vegetation_Description.xlsx"
$ws2.Range("D20:F20").ClearContents()

[void]$ws2.Activate()
[void]$ws2.Range("F14").Select()

Write-Host "Done."
